# Update '想去人数' (F column) values across sheets per the source data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1291
$ws.Range("F6").Value = 321
$ws.Range("F7").Value = 1131
$ws.Range("F8").Value = 434
$ws.Range("F9").Value = 6984
$ws.Range("F11").Value = 87
$ws.Range("F12").Value = 2036
$ws.Range("F13").Value = 7882
$ws.Range("F16").Value = 5461
$ws.Range("F18").Value = 2342
$ws.Range("F19").Value = 994
$ws.Range("F20").Value = 4542
$ws.Range("F21").Value = 278
$ws.Range("F22").Value = 373
$ws.Range("F23").Value = 76
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 333
$ws.Range("F28").Value = 2157
$ws.Range("F30").Value = 243
$ws.Range("F31").Value = 71
$ws.Range("F32").Value = 77
$ws.Range("F33").Value = 553
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 24
$ws.Range("F36").Value = 1435
$ws.Range("F39").Value = 2191
$ws.Range("F40").Value = 2191

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 45
$ws.Range("F6").Value = 24
$ws.Range("F7").Value = 92

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1268

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1268
$ws.Range("F6").Value = 1291
$ws.Range("F9").Value = 321
$ws.Range("F10").Value = 1131
$ws.Range("F11").Value = 434
$ws.Range("F12").Value = 6984
$ws.Range("F14").Value = 87
$ws.Range("F15").Value = 2036
$ws.Range("F16").Value = 7882
$ws.Range("F19").Value = 5461
$ws.Range("F21").Value = 2343
$ws.Range("F22").Value = 994
$ws.Range("F23").Value = 4542
$ws.Range("F24").Value = 278
$ws.Range("F25").Value = 373
$ws.Range("F26").Value = 76
$ws.Range("F28").Value = 6
$ws.Range("F29").Value = 45
$ws.Range("F30").Value = 333
$ws.Range("F33").Value = 2157
$ws.Range("F35").Value = 243
$ws.Range("F36").Value = 71
$ws.Range("F37").Value = 77
$ws.Range("F38").Value = 553
$ws.Range("F39").Value = 1
$ws.Range("F40").Value = 24
$ws.Range("F42").Value = 1435
$ws.Range("F45").Value = 2191
$ws.Range("F46").Value = 24
$ws.Range("F47").Value = 2191
$ws.Range("F49").Value = 92
